$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price (D) column range we are touching so that
# numeric-looking strings (e.g. "16.20", "0.2660") keep their exact text
# representation (including trailing zeros) instead of being coerced into
# floating point numbers by Excel's automatic type inference.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.327.16'
$ws.Range('E2').Value = '  -0.84%  '
$ws.Range('D3').Value = '1.707.51'
$ws.Range('E3').Value = '  -0.93%  '
$ws.Range('D5').Value = '224.02'
$ws.Range('E5').Value = '  -0.85%  '
$ws.Range('D6').Value = '0.5318'
$ws.Range('E6').Value = '  -1.34%  '
$ws.Range('D7').Value = '1.003'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '0.2660'
$ws.Range('E8').Value = '  -0.74%  '
$ws.Range('D9').Value = '0.06597'
$ws.Range('E9').Value = '  -0.31%  '
$ws.Range('D10').Value = '20.76'
$ws.Range('E10').Value = '  -4.71%  '
$ws.Range('D11').Value = '0.07671'
$ws.Range('E11').Value = '  -0.64%  '
$ws.Range('D12').Value = '4.507'
$ws.Range('E12').Value = '  -2.35%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '1.940.53'
$ws.Range('E13').Value = '  -1.05%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.693.95'
$ws.Range('E14').Value = '  -0.88%  '
$ws.Range('D15').Value = '0.5821'
$ws.Range('E15').Value = '  -0.91%  '
$ws.Range('D16').Value = '0.0₅8171'
$ws.Range('E16').Value = '  -1.78%  '
$ws.Range('D17').Value = '67.63'
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('D18').Value = '27.318.93'
$ws.Range('E18').Value = '  -0.94%  '
$ws.Range('D19').Value = '214.94'
$ws.Range('E19').Value = '  -3.09%  '
$ws.Range('D20').Value = '1.003'
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').Value = '4.626'
$ws.Range('E21').Value = '  -2.44%  '
$ws.Range('D22').Value = '10.41'
$ws.Range('E22').Value = '  -2.77%  '
$ws.Range('D23').Value = '5.991'
$ws.Range('E23').Value = '  -1.81%  '
$ws.Range('D24').Value = '1.003'
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('D25').Value = '143.82'
$ws.Range('E25').Value = '  -3.00%  '
$ws.Range('D26').Value = '1.683'
$ws.Range('E26').Value = '  -0.80%  '
$ws.Range('D27').Value = '0.1205'
$ws.Range('E27').Value = '  -2.46%  '
$ws.Range('D28').Value = '7.234'
$ws.Range('E28').Value = '  -2.32%  '
$ws.Range('D29').Value = '16.20'
$ws.Range('E29').Value = '  -2.79%  '
$ws.Range('D30').Value = '0.05368'
$ws.Range('E30').Value = '  -3.52%  '
$ws.Range('D31').Value = '1.291'
$ws.Range('E31').Value = '  -0.99%  '
$ws.Range('D32').Value = '3.487'
$ws.Range('E32').Value = '  -1.83%  '
$ws.Range('D33').Value = '3.420'
$ws.Range('E33').Value = '  -1.13%  '
$ws.Range('D34').Value = '1.647'
$ws.Range('E34').Value = '  -0.88%  '
$ws.Range('D35').Value = '2.861'
$ws.Range('E35').Value = '  +1.41%  '
$ws.Range('D36').Value = '0.9513'
$ws.Range('E36').Value = '  -1.31%  '
$ws.Range('E37').Value = '  -2.04%  '
$ws.Range('D38').Value = '0.5850'
$ws.Range('E38').Value = '  -1.84%  '
$ws.Range('D39').Value = '0.01641'
$ws.Range('E39').Value = '  -0.48%  '
$ws.Range('D40').Value = '5.809'
$ws.Range('E40').Value = '  -2.03%  '
$ws.Range('D41').Value = '1.046.10'
$ws.Range('E41').Value = '  -0.99%  '
$ws.Range('D42').Value = '1.003'
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('D43').Value = '0.8422'
$ws.Range('D44').Value = '100.92'
$ws.Range('E44').Value = '  -0.57%  '
$ws.Range('D45').Value = '1.849.50'
$ws.Range('E45').Value = '  -0.99%  '
$ws.Range('E46').Value = '  +1.92%  '
$ws.Range('D47').Value = '57.86'
$ws.Range('E47').Value = '  -2.15%  '
$ws.Range('D48').Value = '0.4523'
$ws.Range('E48').Value = '  +1.85%  '
$ws.Range('D49').Value = '1.005'
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('D50').Value = '8.097'
$ws.Range('E50').Value = '  -1.15%  '
$ws.Range('E51').Value = '  -0.83%  '

# Restore the Price column to the workbook's default (unstyled) cell style
# now that the text values are safely in place, so no stray number format
# is left applied to the cells.
$ws.Range("D2:D51").Style = "Normal"
